$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.165.96'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.127.63'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E3').Style = "Normal"

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E4').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.48'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '174.41'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('E6').Style = "Normal"

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E7').Style = "Normal"

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('E8').Style = "Normal"

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.156'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('E9').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.482'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('E11').Style = "Normal"

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('E12').Style = "Normal"

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '37.41'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('E13').Style = "Normal"

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -1.64%  '
$ws.Range('E14').Style = "Normal"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.646.11'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('E15').Style = "Normal"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '67.137.93'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('E16').Style = "Normal"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.14'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.89%  '
$ws.Range('E17').Style = "Normal"

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.126.50'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('E18').Style = "Normal"

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '16.48'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('E19').Style = "Normal"

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '493.53'
$ws.Range('D20').Style = "Normal"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.709'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('E21').Style = "Normal"

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +4.80%  '
$ws.Range('E22').Style = "Normal"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.33'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.42%  '
$ws.Range('E24').Style = "Normal"

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -4.11%  '
$ws.Range('E25').Style = "Normal"

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +3.36%  '
$ws.Range('E26').Style = "Normal"

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E27').Style = "Normal"

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('E28').Style = "Normal"

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.66%  '
$ws.Range('E29').Style = "Normal"

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '28.75'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('E31').Style = "Normal"

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.115'
$ws.Range('D32').Style = "Normal"

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0₃0950'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -6.47%  '
$ws.Range('E33').Style = "Normal"

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E34').Style = "Normal"

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.91'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('E35').Style = "Normal"

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.977'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -2.52%  '
$ws.Range('E36').Style = "Normal"

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '47.13'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('E37').Style = "Normal"

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -3.36%  '
$ws.Range('E38').Style = "Normal"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '50.16'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('E39').Style = "Normal"

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('E40').Style = "Normal"

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +1.84%  '
$ws.Range('E41').Style = "Normal"

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.43%  '
$ws.Range('E42').Style = "Normal"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.838.11'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('E43').Style = "Normal"

$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('B44').Style = "Normal"
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C44').Style = "Normal"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.62'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -6.94%  '
$ws.Range('E44').Style = "Normal"

$ws.Range('B45').NumberFormat = "@"
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('B45').Style = "Normal"
$ws.Range('C45').NumberFormat = "@"
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('C45').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '385.43'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('E45').Style = "Normal"

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -2.19%  '
$ws.Range('E46').Style = "Normal"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '135.83'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('E47').Style = "Normal"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '25.04'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('E49').Style = "Normal"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.22'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('E50').Style = "Normal"

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.65%  '
$ws.Range('E51').Style = "Normal"
